$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-29 00:44:29"
$wsZhCn.Range("H4").Value = "2016-08-29 00:44:25"
$wsZhCn.Range("K4").Value = "2016-08-29 00:44:59"
$wsDeDe.Range("K4").Value = "2016-08-29 00:45:12"
